$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the existing row 691 (2019-11-29), shifting
# all subsequent rows down by 9 (old row 691..763 -> new row 700..772).
$ws.Rows("691:699").Insert()

$r = 691
$ws.Cells.Item($r, 1).Value = 1574035200
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-18"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.07000000000000001
$ws.Cells.Item($r, 6).Value = 0.075
$ws.Cells.Item($r, 7).Value = 0.07000000000000001
$ws.Cells.Item($r, 8).Value = 0.07000000000000001
$ws.Cells.Item($r, 9).Value = 2017000

$r = 692
$ws.Cells.Item($r, 1).Value = 1574121600
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-19"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.07000000000000001
$ws.Cells.Item($r, 6).Value = 0.075
$ws.Cells.Item($r, 7).Value = 0.07000000000000001
$ws.Cells.Item($r, 8).Value = 0.07000000000000001
$ws.Cells.Item($r, 9).Value = 300000

$r = 693
$ws.Cells.Item($r, 1).Value = 1574208000
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-20"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.07000000000000001
$ws.Cells.Item($r, 6).Value = 0.08
$ws.Cells.Item($r, 7).Value = 0.07000000000000001
$ws.Cells.Item($r, 8).Value = 0.075
$ws.Cells.Item($r, 9).Value = 14013100

$r = 694
$ws.Cells.Item($r, 1).Value = 1574294400
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-21"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.075
$ws.Cells.Item($r, 6).Value = 0.08
$ws.Cells.Item($r, 7).Value = 0.07000000000000001
$ws.Cells.Item($r, 8).Value = 0.07000000000000001
$ws.Cells.Item($r, 9).Value = 2553900

$r = 695
$ws.Cells.Item($r, 1).Value = 1574380800
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-22"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.075
$ws.Cells.Item($r, 6).Value = 0.08
$ws.Cells.Item($r, 7).Value = 0.075
$ws.Cells.Item($r, 8).Value = 0.08
$ws.Cells.Item($r, 9).Value = 9401100

$r = 696
$ws.Cells.Item($r, 1).Value = 1574640000
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-25"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.08
$ws.Cells.Item($r, 6).Value = 0.08500000000000001
$ws.Cells.Item($r, 7).Value = 0.075
$ws.Cells.Item($r, 8).Value = 0.08500000000000001
$ws.Cells.Item($r, 9).Value = 7416100

$r = 697
$ws.Cells.Item($r, 1).Value = 1574726400
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-26"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.08
$ws.Cells.Item($r, 6).Value = 0.08500000000000001
$ws.Cells.Item($r, 7).Value = 0.08
$ws.Cells.Item($r, 8).Value = 0.08
$ws.Cells.Item($r, 9).Value = 1412400

$r = 698
$ws.Cells.Item($r, 1).Value = 1574812800
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-27"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.08
$ws.Cells.Item($r, 6).Value = 0.08
$ws.Cells.Item($r, 7).Value = 0.075
$ws.Cells.Item($r, 8).Value = 0.075
$ws.Cells.Item($r, 9).Value = 1465000

$r = 699
$ws.Cells.Item($r, 1).Value = 1574899200
$ws.Cells.Item($r, 2).NumberFormat = "@"
$ws.Cells.Item($r, 2).Value = "2019-11-28"
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = "0189"
$ws.Cells.Item($r, 4).Value = "MATANG"
$ws.Cells.Item($r, 5).Value = 0.075
$ws.Cells.Item($r, 6).Value = 0.075
$ws.Cells.Item($r, 7).Value = 0.075
$ws.Cells.Item($r, 8).Value = 0.075
$ws.Cells.Item($r, 9).Value = 383000

# The date strings in column B and the zero-padded id strings in column C
# would otherwise be auto-converted by Excel (dates / numbers). Clear the
# number-format styling applied to force text, leaving plain text cells
# with no explicit style, matching the rest of the sheet.
$ws.Range("B691:C699").ClearFormats()
